$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 1085.5
$ws.Range("J127").Value = 1320.091
$ws.Range("L127").Value = 3960.273
$ws.Range("N127").Value = -13880.273

$ws.Range("H137").Value = 1193.4912
$ws.Range("I137").Value = 1160.909
$ws.Range("J137").Value = 1213.9714
$ws.Range("K137").Value = 3482.727
$ws.Range("L137").Value = 3641.9142
$ws.Range("M137").Value = -932.7270000000003
$ws.Range("N137").Value = -8741.914199999999

$ws.Range("H138").Value = 1963.069
$ws.Range("I138").Value = 1186.9166
$ws.Range("J138").Value = 2918.3333
$ws.Range("K138").Value = 3560.7498
$ws.Range("L138").Value = 8754.999899999999
$ws.Range("M138").Value = 1579.2502
$ws.Range("N138").Value = -19034.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7499.07
$ws.Range("I32").Value = 6022.1445
$ws.Range("J32").Value = 20791.4
$ws.Range("K32").Value = 6022.1445
$ws.Range("L32").Value = 20791.4
$ws.Range("M32").Value = -5735.1445
$ws.Range("N32").Value = -21365.4

$ws.Range("H101").Value = 20833
$ws.Range("J101").Value = 20833
$ws.Range("L101").Value = 20833
$ws.Range("N101").Value = -27323

$ws.Range("H105").Value = 36266.332
$ws.Range("J105").Value = 36266.332
$ws.Range("L105").Value = 36266.332
$ws.Range("N105").Value = -43254.332

$ws.Range("H106").Value = 31666.666
$ws.Range("J106").Value = 31666.666
$ws.Range("L106").Value = 31666.666
$ws.Range("N106").Value = -34190.666

$ws.Range("H122").Value = 1235.6923
$ws.Range("I122").Value = 1013.75
$ws.Range("J122").Value = 1590.8
$ws.Range("K122").Value = 3041.25
$ws.Range("L122").Value = 4772.4
$ws.Range("M122").Value = -591.25
$ws.Range("N122").Value = -9672.4

$ws.Range("H132").Value = 1509325
$ws.Range("I132").Value = 2102.709
$ws.Range("J132").Value = 4469940
$ws.Range("K132").Value = 6308.126999999999
$ws.Range("L132").Value = 13409820
$ws.Range("M132").Value = -3778.126999999999
$ws.Range("N132").Value = -13414880

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H109").Value = 29463.334
$ws.Range("J109").Value = 29463.334
$ws.Range("L109").Value = 29463.334
$ws.Range("N109").Value = -32237.334

$ws.Range("H122").Value = 38870
$ws.Range("J122").Value = 38870
$ws.Range("L122").Value = 38870
$ws.Range("N122").Value = -48670

$ws.Range("H134").Value = 2241.037
$ws.Range("I134").Value = 1193.6154
$ws.Range("J134").Value = 3213.6428
$ws.Range("K134").Value = 3580.8462
$ws.Range("L134").Value = 9640.928400000001
$ws.Range("M134").Value = -1045.8462
$ws.Range("N134").Value = -14710.9284

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 5000395
$ws.Range("I6").Value = 10000000
$ws.Range("K6").Value = 10000000
$ws.Range("M6").Value = -9999887

$ws.Range("H28").Value = 25985.8
$ws.Range("J28").Value = 25985.8
$ws.Range("L28").Value = 25985.8
$ws.Range("N28").Value = -26475.8

$ws.Range("H43").Value = 17551.625
$ws.Range("J43").Value = 17551.625
$ws.Range("L43").Value = 17551.625
$ws.Range("N43").Value = -17919.625

$ws.Range("H51").Value = 24772.5
$ws.Range("I51").Value = 10090
$ws.Range("J51").Value = 29666.666
$ws.Range("K51").Value = 10090
$ws.Range("L51").Value = 29666.666
$ws.Range("M51").Value = -9354
$ws.Range("N51").Value = -31138.666

$ws.Range("H61").Value = 24772.5
$ws.Range("I61").Value = 10090
$ws.Range("J61").Value = 29666.666
$ws.Range("K61").Value = 10090
$ws.Range("L61").Value = 29666.666
$ws.Range("M61").Value = -9742
$ws.Range("N61").Value = -30362.666

$ws.Range("H99").Value = 3309.4
$ws.Range("I99").Value = 3525
$ws.Range("J99").Value = 3165.6667
$ws.Range("K99").Value = 3525
$ws.Range("L99").Value = 3165.6667
$ws.Range("M99").Value = -2027
$ws.Range("N99").Value = -6161.6667

$ws.Range("H101").Value = 17551.625
$ws.Range("J101").Value = 17551.625
$ws.Range("L101").Value = 17551.625
$ws.Range("N101").Value = -24041.625

$ws.Range("H126").Value = 3309.4
$ws.Range("I126").Value = 3525
$ws.Range("J126").Value = 3165.6667
$ws.Range("K126").Value = 10575
$ws.Range("L126").Value = 9497.000100000001
$ws.Range("M126").Value = -8105
$ws.Range("N126").Value = -14437.0001

$ws.Range("H132").Value = 1946.25
$ws.Range("I132").Value = 1296.7142
$ws.Range("J132").Value = 2855.6
$ws.Range("K132").Value = 3890.1426
$ws.Range("L132").Value = 8566.799999999999
$ws.Range("M132").Value = -1360.1426
$ws.Range("N132").Value = -13626.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 295.63635
$ws.Range("J7").Value = 324.2
$ws.Range("L7").Value = 972.5999999999999
$ws.Range("N7").Value = -1196.6

$ws.Range("H131").Value = 338549.44
$ws.Range("I131").Value = 429.0909
$ws.Range("J131").Value = 834459.25
$ws.Range("K131").Value = 1287.2727
$ws.Range("L131").Value = 2503377.75
$ws.Range("M131").Value = 3752.7273
$ws.Range("N131").Value = -2513457.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 16654.166
$ws.Range("J123").Value = 16654.166
$ws.Range("L123").Value = 16654.166
$ws.Range("N123").Value = -21554.166

$ws.Range("H132").Value = 2586.2554
$ws.Range("I132").Value = 1513.375
$ws.Range("J132").Value = 3705.7827
$ws.Range("K132").Value = 4540.125
$ws.Range("L132").Value = 11117.3481
$ws.Range("M132").Value = -2010.125
$ws.Range("N132").Value = -16177.3481

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 24999.666
$ws.Range("J64").Value = 24999.666
$ws.Range("L64").Value = 24999.666
$ws.Range("N64").Value = -25449.666

$ws.Range("H67").Value = 24999.666
$ws.Range("J67").Value = 24999.666
$ws.Range("L67").Value = 24999.666
$ws.Range("N67").Value = -26559.666

$ws.Range("H111").Value = 29546.75
$ws.Range("J111").Value = 29546.75
$ws.Range("L111").Value = 29546.75
$ws.Range("N111").Value = -37726.75

$ws.Range("H132").Value = 19763.59
$ws.Range("I132").Value = 27052.762
$ws.Range("J132").Value = 3650.6843
$ws.Range("K132").Value = 81158.28599999999
$ws.Range("L132").Value = 10952.0529
$ws.Range("M132").Value = -78628.28599999999
$ws.Range("N132").Value = -16012.0529

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 26299.8
$ws.Range("J63").Value = 26299.8
$ws.Range("L63").Value = 26299.8
$ws.Range("N63").Value = -27547.8

$ws.Range("H66").Value = 26299.8
$ws.Range("J66").Value = 26299.8
$ws.Range("L66").Value = 78899.39999999999
$ws.Range("N66").Value = -85139.39999999999

$ws.Range("H132").Value = 1754.421
$ws.Range("I132").Value = 1480
$ws.Range("J132").Value = 2131.75
$ws.Range("K132").Value = 4440
$ws.Range("L132").Value = 6395.25
$ws.Range("M132").Value = -1910
$ws.Range("N132").Value = -11455.25

